$wb = $excel.ActiveWorkbook

# --- "System Settings" sheet: insert a new row at row 9 (TFM_MIG table gains a row) ---
$ws = $wb.Worksheets.Item("System Settings")

$ws.Rows.Item(9).Insert()

# Copy the formatting of the row above (row 8, the last existing data row in
# this small table) onto only the relevant cells of the freshly inserted row,
# so the new cells carry the same style without touching the whole row.
$ws.Range("C8:E8").Copy()
$ws.Range("C9:E9").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(9, 3).Value = "NCAP_AF"
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 3

# Make "System Settings" the active sheet / tab, with C11 selected.
$ws.Activate()
$ws.Range("C11").Select()

# --- "fuels" sheet: no longer the active tab ---
$fuels = $wb.Worksheets.Item("fuels")
$fuels.Range("F20").Select()

# Re-activate "System Settings" last so it ends up as the workbook's active sheet.
$ws.Activate()
